{"js": "// Office.js (Word JavaScript API) script.\n// Implements the same content edit described by the XML diff:\n//   1. \"Lecture 8: R basics, ...\"  -> \"Lecture 8: R and RStudio basics, ...\"\n//   2. \"Lecture 9: The RStudio environment, ...\" -> \"Lecture 10: Accessing help, ...\"\n//      (keeps its own trailing \" pdf\" hyperlink run untouched)\n//   3. The paragraph that used to read \"Lecture 10: Accessing help, ...\" (with its\n//      own \" pdf\" hyperlink) is removed entirely, since that text now lives in the\n//      paragraph from step 2.\n\nconst body = context.document.body;\n\n// 1) Lecture 8 text update (only the bold \"title\" run; the trailing \" \" + hyperlink\n//    run are untouched, exactly like the diff).\nconst lecture8 = body.search(\n  \"Lecture 8: R basics, Brian, XX min (much live demo)\",\n  { matchCase: true }\n);\nlecture8.load(\"items\");\nawait context.sync();\n\nif (lecture8.items.length === 0) {\n  throw new Error(\"Could not find Lecture 8 heading text to update\");\n}\nlecture8.items[0].insertText(\n  \"Lecture 8: R and RStudio basics, Brian, XX min (much live demo)\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 2) Lecture 9 paragraph's bold title run becomes the old Lecture 10 title text.\nconst lecture9 = body.search(\n  \"Lecture 9: The RStudio environment, Brian, XX min (much live demo)\",\n  { matchCase: true }\n);\nlecture9.load(\"items\");\nawait context.sync();\n\nif (lecture9.items.length === 0) {\n  throw new Error(\"Could not find Lecture 9 heading text to update\");\n}\nlecture9.items[0].insertText(\n  \"Lecture 10: Accessing help, Brian, XX min (much live demo)\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 3) Remove the paragraph that used to hold \"Lecture 10: Accessing help, ...\" text\n//    (and its own \"pdf\" hyperlink). After step 2 there are two paragraphs that start\n//    with that text; delete the later (original) one, keep the renamed-in-place one.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst matchingIndexes = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Lecture 10: Accessing help\") === 0) {\n    matchingIndexes.push(i);\n  }\n}\n\nif (matchingIndexes.length < 2) {\n  throw new Error(\n    \"Expected two paragraphs starting with 'Lecture 10: Accessing help' before cleanup, found \" +\n      matchingIndexes.length\n  );\n}\n\n// Delete the last (originally-existing) \"Lecture 10\" paragraph -- this is the one\n// that the diff removes outright.\nconst deleteIndex = matchingIndexes[matchingIndexes.length - 1];\nparagraphs.items[deleteIndex].delete();\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Implements the same content edit described by the XML diff:\n#   1. \"Lecture 8: R basics, ...\"  -> \"Lecture 8: R and RStudio basics, ...\"\n#   2. \"Lecture 9: The RStudio environment, ...\" -> \"Lecture 10: Accessing help, ...\"\n#      (keeps its own trailing \" pdf\" hyperlink run untouched)\n#   3. The paragraph that used to read \"Lecture 10: Accessing help, ...\" (with its\n#      own \" pdf\" hyperlink) is removed entirely, since that text now lives in the\n#      paragraph from step 2.\n\n$d = $word.ActiveDocument\n\n# 1) Lecture 8 text update (only the bold \"title\" run; the trailing \" \" + hyperlink\n#    run are untouched, exactly like the diff). wdReplaceOne = 1 ... but since there\n#    is only a single match, wdReplaceAll (2) is equally safe and avoids any\n#    ambiguity about current selection/start position.\n$found8 = $d.Content.Find.Execute(\n    \"Lecture 8: R basics, Brian, XX min (much live demo)\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Lecture 8: R and RStudio basics, Brian, XX min (much live demo)\",\n    2\n)\nif (-not $found8) {\n    throw \"Could not find Lecture 8 heading text to update\"\n}\n\n# 2) Lecture 9 paragraph's bold title run becomes the old Lecture 10 title text.\n$found9 = $d.Content.Find.Execute(\n    \"Lecture 9: The RStudio environment, Brian, XX min (much live demo)\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Lecture 10: Accessing help, Brian, XX min (much live demo)\",\n    2\n)\nif (-not $found9) {\n    throw \"Could not find Lecture 9 heading text to update\"\n}\n\n# 3) Remove the paragraph that used to hold \"Lecture 10: Accessing help, ...\" text\n#    (and its own \"pdf\" hyperlink). After step 2 there are two paragraphs that start\n#    with that text; walk all paragraphs and keep the *last* match (the original,\n#    now-duplicate paragraph) to delete, leaving the renamed-in-place one from step 2.\n$targetPara = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text.StartsWith(\"Lecture 10: Accessing help\")) {\n        $targetPara = $para\n    }\n}\n\nif ($targetPara -eq $null) {\n    throw \"Expected a duplicated 'Lecture 10: Accessing help' paragraph to remove\"\n}\n\n$targetPara.Range.Delete()\n"}
